# Applies the NATMI TPM recompute for the Snca-Lag3 sheet (rows 2-13).
# Only the numeric columns E:J (ligand stats) and M:T (receptor/edge stats) change;
# A:D (cluster/ligand/receptor labels) and K:L stay as-is.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=11.1710695; "N"=22.342139; "O"=0.3825233089595566; "P"=0.3232015958207352; "Q"=0.3614567102767501; "R"=1.445826841107; "S"=0.02931007110497226; "T"=0.02476466540173064 }
    3 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=4.147608666666667; "N"=12.442826; "O"=0.1420237329507858; "P"=0.1799980395663877; "Q"=0.134202099823; "R"=0.8052125989380001; "S"=0.01088227988695256; "T"=0.0137919839520269 }
    4 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=2.534243; "N"=7.602729; "O"=0.08677835350210591; "P"=0.1099811502109347; "Q"=0.08199923362950001; "R"=0.4919954017770001; "S"=0.006649214967938227; "T"=0.008427082108164941 }
    5 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=7.312259; "N"=14.624518; "O"=0.2503886945336163; "P"=0.211558416842231; "Q"=0.2365991083335; "R"=0.9463964333340001; "S"=0.01918552482624634; "T"=0.01621023371717394 }
    6 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=2.016508333333333; "N"=6.049525; "O"=0.06904991865024089; "P"=0.08751248633613068; "Q"=0.06524715188750001; "R"=0.391482911325; "S"=0.005290809678855645; "T"=0.006705466404286738 }
    7 = @{ "E"=1; "F"=0.5; "G"=0.0323565; "H"=0.064713; "I"=0.076622967590378; "J"=0.076622967590378; "M"=2.021942333333333; "N"=6.065827; "O"=0.06923599140369445; "P"=0.08774831122358079; "Q"=0.0654229771085; "R"=0.392537862651; "S"=0.005305067125412969; "T"=0.006723536006994832 }
    8 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=11.1710695; "N"=22.342139; "O"=0.3825233089595566; "P"=0.3232015958207352; "Q"=4.35588486032225; "R"=17.423539441289; "S"=0.3532132378545844; "T"=0.2984369304190045 }
    9 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=4.147608666666667; "N"=12.442826; "O"=0.1420237329507858; "P"=0.1799980395663877; "Q"=1.617258383154333; "R"=9.703550298926; "S"=0.1311414530638332; "T"=0.1662060556143608 }
    10 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=2.534243; "N"=7.602729; "O"=0.08677835350210591; "P"=0.1099811502109347; "Q"=0.9881659688965; "R"=5.928995813379; "S"=0.08012913853416768; "T"=0.1015540681027697 }
    11 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=7.312259; "N"=14.624518; "O"=0.2503886945336163; "P"=0.211558416842231; "Q"=2.8512362467045; "R"=11.404944986818; "S"=0.23120316970737; "T"=0.195348183125057 }
    12 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=2.016508333333333; "N"=6.049525; "O"=0.06904991865024089; "P"=0.08751248633613068; "Q"=0.7862880201291667; "R"=4.717728120775; "S"=0.06375910897138524; "T"=0.08080701993184394 }
    13 = @{ "E"=2; "F"=1; "G"=0.3899255; "H"=0.779851; "I"=0.923377032409622; "J"=0.923377032409622; "M"=2.021942333333333; "N"=6.065827; "O"=0.06923599140369445; "P"=0.08774831122358079; "Q"=0.7884068752961665; "R"=4.730441251776999; "S"=0.06393092427828148; "T"=0.08102477521658595 }
}

foreach ($row in $newValues.Keys) {
    foreach ($col in $newValues[$row].Keys) {
        $ws.Range("$col$row").Value = $newValues[$row][$col]
    }
}